$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '243.25'
$c.Style = "Normal"
$c = $ws.Cells.Item(2, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '23.01'
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = '5.397'
$c.Style = "Normal"
$c = $ws.Cells.Item(4, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '0.05911'
$c.Style = "Normal"
$c = $ws.Cells.Item(5, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '3.457'
$c.Style = "Normal"
$c = $ws.Cells.Item(6, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '6.553'
$c.Style = "Normal"
$c = $ws.Cells.Item(7, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(8, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '0.9196'
$c.Style = "Normal"
$c = $ws.Cells.Item(9, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.1411'
$c.Style = "Normal"
$c = $ws.Cells.Item(10, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.07410'
$c.Style = "Normal"
$c = $ws.Cells.Item(11, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.03290'
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '0.03064'
$c.Style = "Normal"
$c = $ws.Cells.Item(13, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = '0.09351'
$c.Style = "Normal"
$c = $ws.Cells.Item(14, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '3.853'
$c.Style = "Normal"
$c = $ws.Cells.Item(15, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '0.001568'
$c.Style = "Normal"
$c = $ws.Cells.Item(16, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.04668'
$c.Style = "Normal"
$c = $ws.Cells.Item(17, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '0.0005953'
$c.Style = "Normal"
$c = $ws.Cells.Item(18, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '0.005963'
$c.Style = "Normal"
$c = $ws.Cells.Item(19, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$ws.Cells.Item(20, 2).Value = 'BitKan'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '0.001319'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '19BitKanKAN'
$c = $ws.Cells.Item(20, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$ws.Cells.Item(21, 2).Value = 'HotbitToken'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '0.004933'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '20HotbitTokenHTB'
$c = $ws.Cells.Item(21, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(22, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '3.618'
$c.Style = "Normal"
$c = $ws.Cells.Item(23, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '2.152'
$c.Style = "Normal"
$c = $ws.Cells.Item(24, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(25, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '0.1327'
$c.Style = "Normal"
$c = $ws.Cells.Item(26, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(27, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(28, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(29, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(30, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(31, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(32, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(33, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(34, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(35, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(36, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(37, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(38, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(39, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.03963'
$c.Style = "Normal"
$c = $ws.Cells.Item(40, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.006207'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '40KickTokenKICKBestin24h'
$c = $ws.Cells.Item(41, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.1074'
$c.Style = "Normal"
$c = $ws.Cells.Item(42, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '0.003001'
$c.Style = "Normal"
$c = $ws.Cells.Item(43, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '0.008103'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '43LocalTradersLCT'
$c = $ws.Cells.Item(44, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.00005190'
$c.Style = "Normal"
$c = $ws.Cells.Item(45, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(46, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '0.7003'
$c.Style = "Normal"
$c = $ws.Cells.Item(47, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '0.002282'
$c.Style = "Normal"
$c = $ws.Cells.Item(48, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '0.00002101'
$c.Style = "Normal"
$c = $ws.Cells.Item(49, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.0002001'
$c.Style = "Normal"
$c = $ws.Cells.Item(50, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"

$c = $ws.Cells.Item(51, 7)
$c.NumberFormat = "@"
$c.Value = '13'
$c.Style = "Normal"
